# "Added results for partial model"
#
# Sheet1 layout (row 1 header): A=Point, B=IPC RO, C=IPC PO, D=DELTA, E=DELTA^2
#   DELTA   = IPC PO (C) - IPC RO (B)
#   DELTA^2 = DELTA squared
#   Row 52  = TOTAL  (sum of DELTA in C52, sum of DELTA^2 in E52)
#   Row 53  = MSE    (mean of DELTA^2 in E53)
#
# The partial model produced new "IPC PO" predictions (column C); this updates
# column C together with the dependent DELTA / DELTA^2 columns and the TOTAL /
# MSE summary rows with the recomputed results. Each triple below is
# (row, column, new value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 2, 29.53999999999999),
    @(2, 3, 29.42570304870605),
    @(2, 4, -0.1142969512939374),
    @(2, 5, 0.01306379307508869),
    @(3, 3, 29.50772666931152),
    @(3, 4, -0.04227333068847372),
    @(3, 5, 0.001787034487497054),
    @(4, 3, 29.72051811218262),
    @(4, 4, -0.02948188781738281),
    @(4, 5, 0.0008691817092767451),
    @(5, 3, 29.87379264831543),
    @(5, 4, 0.03379264831542628),
    @(5, 5, 0.001141943080170082),
    @(6, 3, 29.8944149017334),
    @(6, 4, 0.08441490173339616),
    @(6, 5, 0.007125875634658931),
    @(7, 3, 29.93996238708496),
    @(7, 4, 0.01996238708495923),
    @(7, 5, 0.0003984968981297471),
    @(8, 3, 30.07316398620605),
    @(8, 4, 0.09316398620605071),
    @(8, 5, 0.008679528325801206),
    @(9, 2, 30.03999999999999),
    @(9, 3, 30.21114921569824),
    @(9, 4, 0.1711492156982501),
    @(9, 5, 0.02929205403412615),
    @(10, 2, 30.21000000000001),
    @(10, 3, 30.26265525817871),
    @(10, 4, 0.05265525817870298),
    @(10, 5, 0.002772576213865867),
    @(11, 3, 30.34427452087402),
    @(11, 4, 0.1242745208740246),
    @(11, 5, 0.01544415653846837),
    @(12, 3, 30.36523246765137),
    @(12, 4, -0.01476753234862827),
    @(12, 5, 0.0002180800116677822),
    @(13, 3, 30.50981712341309),
    @(13, 4, 0.06981712341308821),
    @(13, 5, 0.00487443072167839),
    @(14, 3, 30.43168449401855),
    @(14, 4, -0.04831550598144929),
    @(14, 5, 0.002334388118243462),
    @(15, 3, 30.43945121765137),
    @(15, 4, -0.2505487823486305),
    @(15, 5, 0.06277469233638144),
    @(16, 3, 30.4976634979248),
    @(16, 4, -0.2523365020751953),
    @(16, 5, 0.06367371027954505),
    @(17, 3, 30.65968132019043),
    @(17, 4, -0.280318679809568),
    @(17, 5, 0.07857856225017913),
    @(18, 3, 30.76798057556152),
    @(18, 4, -0.1820194244384794),
    @(18, 5, 0.03313107087291531),
    @(19, 3, 31.12988090515137),
    @(19, 4, 0.1098809051513712),
    @(19, 5, 0.01207381331688463),
    @(20, 3, 31.24157524108887),
    @(20, 4, 0.1215752410888626),
    @(20, 5, 0.01478053924581508),
    @(21, 3, 31.3477611541748),
    @(21, 4, 0.06776115417480355),
    @(21, 5, 0.004591574015101497),
    @(22, 3, 31.23878288269043),
    @(22, 4, -0.1412171173095658),
    @(22, 5, 0.01994227422122366),
    @(23, 3, 31.39718437194824),
    @(23, 4, -0.1828156280517561),
    @(23, 5, 0.03342155385995804),
    @(24, 2, 31.65000000000001),
    @(24, 3, 31.92793846130371),
    @(24, 4, 0.2779384613037053),
    @(24, 5, 0.07724978827187126),
    @(25, 3, 32.41188430786133),
    @(25, 4, 0.5318843078613327),
    @(25, 5, 0.2829009169491289),
    @(26, 3, 32.36374282836914),
    @(26, 4, 0.08374282836913949),
    @(26, 5, 0.007012861303263154),
    @(27, 3, 32.45510482788086),
    @(27, 4, 0.005104827880856533),
    @(27, 5, 0.0000260592676931702),
    @(28, 2, 32.84999999999999),
    @(28, 3, 32.71161270141602),
    @(28, 4, -0.1383872985839787),
    @(28, 5, 0.01915104440937127),
    @(29, 2, 32.90000000000001),
    @(29, 3, 32.94222640991211),
    @(29, 4, 0.04222640991210369),
    @(29, 5, 0.001783069694065009),
    @(30, 2, 33.09999999999999),
    @(30, 3, 32.91791915893555),
    @(30, 4, -0.1820808410644474),
    @(30, 5, 0.03315343268273657),
    @(31, 2, 33.40000000000001),
    @(31, 3, 33.66357040405273),
    @(31, 4, 0.2635704040527287),
    @(31, 5, 0.06946935789251867),
    @(32, 3, 33.69541549682617),
    @(32, 4, -0.004584503173830967),
    @(32, 5, 0.00002101766935086621),
    @(33, 2, 34.09999999999999),
    @(33, 3, 33.89573287963867),
    @(33, 4, -0.2042671203613224),
    @(33, 5, 0.04172505646070699),
    @(34, 2, 34.40000000000001),
    @(34, 3, 34.42705917358398),
    @(34, 4, 0.02705917358397869),
    @(34, 5, 0.0007321988750478901),
    @(35, 2, 34.90000000000001),
    @(35, 3, 35.06875228881836),
    @(35, 4, 0.1687522888183537),
    @(35, 5, 0.02847733498143306),
    @(36, 3, 35.66774368286133),
    @(36, 4, 0.367743682861331),
    @(36, 5, 0.1352354162844152),
    @(37, 3, 35.96606826782227),
    @(37, 4, 0.2660682678222628),
    @(37, 5, 0.07079232314193935),
    @(38, 3, 35.86572647094727),
    @(38, 4, -0.4342735290527315),
    @(38, 5, 0.1885934980359137),
    @(39, 3, 36.4784049987793),
    @(39, 4, -0.3215950012207003),
    @(39, 5, 0.1034233448101422),
    @(40, 3, 37.2089958190918),
    @(40, 4, -0.09100418090820028),
    @(40, 5, 0.008281760942772444),
    @(41, 2, 37.90000000000001),
    @(41, 3, 37.99245071411133),
    @(41, 4, 0.09245071411132244),
    @(41, 5, 0.008547134539693474),
    @(42, 3, 38.43264389038086),
    @(42, 4, -0.06735610961914062),
    @(42, 5, 0.004536845503025688),
    @(43, 2, 38.90000000000001),
    @(43, 3, 39.07284927368164),
    @(43, 4, 0.1728492736816349),
    @(43, 5, 0.02987687141226874),
    @(44, 2, 39.40000000000001),
    @(44, 3, 39.60612869262695),
    @(44, 4, 0.2061286926269474),
    @(44, 5, 0.04248903792409458),
    @(45, 2, 39.90000000000001),
    @(45, 3, 39.64629745483398),
    @(45, 4, -0.2537025451660213),
    @(45, 5, 0.06436498142371708),
    @(46, 2, 40.09999999999999),
    @(46, 3, 39.90364456176758),
    @(46, 4, -0.1963554382324162),
    @(46, 5, 0.03855545812344421),
    @(47, 2, 40.59999999999999),
    @(47, 3, 40.41103744506836),
    @(47, 4, -0.1889625549316349),
    @(47, 5, 0.03570684716629115),
    @(48, 2, 40.90000000000001),
    @(48, 3, 40.64894485473633),
    @(48, 4, -0.2510551452636776),
    @(48, 5, 0.06302868596336623),
    @(49, 2, 41.20000000000001),
    @(49, 3, 41.17435073852539),
    @(49, 4, -0.02564926147461932),
    @(49, 5, 0.000657884614193391),
    @(50, 3, 41.48037338256836),
    @(50, 4, -0.01962661743164062),
    @(50, 5, 0.0003852041118079796),
    @(51, 3, 42.27613067626953),
    @(51, 4, 0.4761306762695341),
    @(51, 5, 0.2267004208848839),
    @(52, 3, 0.0128058624267382),
    @(52, 5, 1.993847182585832),
    @(53, 5, 0.03987694365171665)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
